$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update surnames (B column) and emails (C column) while keeping first
# names (A column) unchanged. Wiktor and Maciej get the male surname
# "Kowalski", Luiza gets the female surname "Kowalska".
$ws.Range("B1").Value = "Kowalski"
$ws.Range("C1").Value = "wiktorkowalski1@gmail.com"

$ws.Range("B2").Value = "Kowalski"
$ws.Range("C2").Value = "maciejkowalski1@gmail.com"

$ws.Range("B3").Value = "Kowalska"
$ws.Range("C3").Value = "luizakowalska1@gmail.com"

# Move the active cell selection from D16 to D6.
$ws.Range("D6").Select() | Out-Null
